$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# 1) sim_ts (D7): 1 -> 5
$ws.Range("D7").Value = 5

# 2) EV_present (D29): 0 -> 40
$ws.Range("D29").Value = 40

# 3) Remove the EV_statut (row31) and EV_nb_drivers (row32) rows entirely.
#    This shifts EV_charger_power / EV_usage / blank / Plot / plt_plot rows up by two.
$ws.Rows("31:32").Delete()

# 4) EV_charger_power description gains a unit suffix.
$ws.Range("B31").Value = "Puissance de charge du chargeur [kW]"

# 5) EV_usage gains a helper note describing accepted values.
$ws.Range("G32").Value = "[short, normal, long, int: (km/year)]"
